# Add a "Color" column (E) to the Gantt timeline sheet, tag a task row with
# a color code, and correct the "Design and train GNN" start date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Color" header in E1, styled like the other table headings (bold
# heading font) but without the heavy bottom border the other headers have.
$ws.Range("E1").Value = "Color"
$ws.Range("E1").Style = "Heading 2"
$ws.Range("E1").Borders.LineStyle = -4142

# Tag the "Implement Multi-vehicle routing solver for GT" task (row 4) with
# its timeline bar color.
$ws.Range("E4").Value = "#e2684d"

# Correct the start date for "Design and train GNN" (row 5) -> 2022-11-15.
$ws.Range("B5").Value = 44880

# Leave the active selection on the newly edited cell, as in the source file.
$ws.Range("E4").Select() | Out-Null
